# Scheduled runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# figures across the job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H15").Value = 1033598.56
$ws.Range("I15").Value = 1033598.56
$ws.Range("K15").Value = 3100795.68
$ws.Range("M15").Value = -3100626.68

$ws.Range("H17").Value = 962.46875
$ws.Range("J17").Value = 962.46875
$ws.Range("L17").Value = 2887.40625
$ws.Range("N17").Value = -3223.40625

$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52496

$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -162480

$ws.Range("H112").Value = 2518.5
$ws.Range("I112").Value = 792.5
$ws.Range("J112").Value = 2710.2778
$ws.Range("K112").Value = 2377.5
$ws.Range("L112").Value = 8130.8334
$ws.Range("M112").Value = -1269.5
$ws.Range("N112").Value = -10346.8334

$ws.Range("H113").Value = 3014.2307
$ws.Range("I113").Value = 2290.625
$ws.Range("J113").Value = 4172
$ws.Range("K113").Value = 2290.625
$ws.Range("L113").Value = 4172
$ws.Range("M113").Value = 963.375
$ws.Range("N113").Value = -10680

$ws.Range("H116").Value = 2833.4119
$ws.Range("J116").Value = 2972.75
$ws.Range("L116").Value = 2972.75
$ws.Range("N116").Value = -9856.75

$ws.Range("H120").Value = 36761
$ws.Range("J120").Value = 36761
$ws.Range("L120").Value = 36761
$ws.Range("N120").Value = -46437

$ws.Range("H125").Value = 658.9091
$ws.Range("I125").Value = 649.7778
$ws.Range("J125").Value = 700
$ws.Range("K125").Value = 5848.000199999999
$ws.Range("L125").Value = 6300
$ws.Range("M125").Value = -3388.000199999999
$ws.Range("N125").Value = -11220

$ws.Range("H132").Value = 6083.769
$ws.Range("I132").Value = 5121.2705
$ws.Range("K132").Value = 15363.8115
$ws.Range("M132").Value = -12833.8115

# ---------------------------------------------------------------------------
# ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()

$ws.Range("H61").Value = 2534.6667
$ws.Range("I61").Value = 4032.625
$ws.Range("J61").Value = 1903.9474
$ws.Range("K61").Value = 4032.625
$ws.Range("L61").Value = 1903.9474
$ws.Range("M61").Value = -3820.625
$ws.Range("N61").Value = -2327.9474

$ws.Range("H96").Value = 25246.4
$ws.Range("J96").Value = 25246.4
$ws.Range("L96").Value = 25246.4
$ws.Range("N96").Value = -30738.4

$ws.Range("H136").Value = 2534.6667
$ws.Range("I136").Value = 4032.625
$ws.Range("J136").Value = 1903.9474
$ws.Range("K136").Value = 12097.875
$ws.Range("L136").Value = 5711.8422
$ws.Range("M136").Value = -9547.875
$ws.Range("N136").Value = -10811.8422

# ---------------------------------------------------------------------------
# BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H63").Value = 50000
$ws.Range("J63").Value = 50000
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51372

$ws.Range("H66").Value = 50000
$ws.Range("J66").Value = 50000
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -156864

$ws.Range("H102").Value = 9700.6
$ws.Range("I102").Value = 6375.75
$ws.Range("J102").Value = 23000
$ws.Range("K102").Value = 6375.75
$ws.Range("L102").Value = 23000
$ws.Range("M102").Value = -3130.75
$ws.Range("N102").Value = -29490

# ---------------------------------------------------------------------------
# CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H4").Value = 10572.875
$ws.Range("J4").Value = 10572.875
$ws.Range("L4").Value = 10572.875
$ws.Range("N4").Value = -10796.875

$ws.Range("H43").Value = 24200
$ws.Range("J43").Value = 24200
$ws.Range("L43").Value = 24200
$ws.Range("N43").Value = -24568

$ws.Range("H58").Value = 6731.136
$ws.Range("I58").Value = 3391.25
$ws.Range("J58").Value = 10739
$ws.Range("K58").Value = 3391.25
$ws.Range("L58").Value = 10739
$ws.Range("M58").Value = -3188.25
$ws.Range("N58").Value = -11145

$ws.Range("H101").Value = 24200
$ws.Range("J101").Value = 24200
$ws.Range("L101").Value = 24200
$ws.Range("N101").Value = -30690

$ws.Range("H134").Value = 2398.3684
$ws.Range("I134").Value = 1541.2858
$ws.Range("J134").Value = 2898.3333
$ws.Range("K134").Value = 4623.857400000001
$ws.Range("L134").Value = 8694.999899999999
$ws.Range("M134").Value = -2088.857400000001
$ws.Range("N134").Value = -13764.9999

$ws.Range("H136").Value = 6731.136
$ws.Range("I136").Value = 3391.25
$ws.Range("J136").Value = 10739
$ws.Range("K136").Value = 10173.75
$ws.Range("L136").Value = 32217
$ws.Range("M136").Value = -7623.75
$ws.Range("N136").Value = -37317

# ---------------------------------------------------------------------------
# CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H59").Value = 1965.6666
$ws.Range("I59").Value = 1400
$ws.Range("J59").Value = 2248.5
$ws.Range("K59").Value = 4200
$ws.Range("L59").Value = 6745.5
$ws.Range("M59").Value = -3660
$ws.Range("N59").Value = -7825.5

$ws.Range("H63").Value = 4971.2856
$ws.Range("I63").Value = 2933
$ws.Range("J63").Value = 6500
$ws.Range("K63").Value = 8799
$ws.Range("L63").Value = 19500
$ws.Range("M63").Value = -8050
$ws.Range("N63").Value = -20998

$ws.Range("H66").Value = 4971.2856
$ws.Range("I66").Value = 2933
$ws.Range("J66").Value = 6500
$ws.Range("K66").Value = 26397
$ws.Range("L66").Value = 58500
$ws.Range("M66").Value = -22653
$ws.Range("N66").Value = -65988

$ws.Range("H75").Value = 4542.6
$ws.Range("J75").Value = 7166.6665
$ws.Range("L75").Value = 21499.9995
$ws.Range("N75").Value = -23495.9995

$ws.Range("H78").Value = 4542.6
$ws.Range("J78").Value = 7166.6665
$ws.Range("L78").Value = 64499.9985
$ws.Range("N78").Value = -74483.9985

$ws.Range("H80").Value = 3214.7144
$ws.Range("J80").Value = 3214.7144
$ws.Range("L80").Value = 9644.143199999999
$ws.Range("N80").Value = -11516.1432

$ws.Range("H83").Value = 3214.7144
$ws.Range("J83").Value = 3214.7144
$ws.Range("L83").Value = 28932.4296
$ws.Range("N83").Value = -38292.4296

$ws.Range("H105").Value = 11115.75
$ws.Range("J105").Value = 11985.714
$ws.Range("L105").Value = 35957.142
$ws.Range("N105").Value = -41199.142

# ---------------------------------------------------------------------------
# GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H5").Value = 10005
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 10005
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 10005
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -10229

$ws.Range("H97").Value = 499.4
$ws.Range("I97").Value = 374.25
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 374.25
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = 121.75
$ws.Range("N97").Value = -1992

$ws.Range("H102").Value = 1751418.2
$ws.Range("I102").Value = 2389068.5
$ws.Range("J102").Value = 12371.637
$ws.Range("K102").Value = 2389068.5
$ws.Range("L102").Value = 12371.637
$ws.Range("M102").Value = -2387446.5
$ws.Range("N102").Value = -15615.637

# ---------------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H2").Value = 2731724.2
$ws.Range("I2").Value = 1001
$ws.Range("J2").Value = 3004796.8
$ws.Range("K2").Value = 1001
$ws.Range("L2").Value = 3004796.8
$ws.Range("M2").Value = -889
$ws.Range("N2").Value = -3005020.8

$ws.Range("H40").Value = 71432024
$ws.Range("I40").Value = 111113736
$ws.Range("J40").Value = 4944.8
$ws.Range("K40").Value = 111113736
$ws.Range("L40").Value = 4944.8
$ws.Range("M40").Value = -111113600
$ws.Range("N40").Value = -5216.8

$ws.Range("H122").Value = 9232.777
$ws.Range("I122").Value = 11409.25
$ws.Range("J122").Value = 4879.8335
$ws.Range("K122").Value = 34227.75
$ws.Range("L122").Value = 14639.5005
$ws.Range("M122").Value = -31777.75
$ws.Range("N122").Value = -19539.5005

$ws.Range("H132").Value = 30654.703
$ws.Range("I132").Value = 33573.516
$ws.Range("J132").Value = 6574.5
$ws.Range("K132").Value = 100720.548
$ws.Range("L132").Value = 19723.5
$ws.Range("M132").Value = -98190.54800000001
$ws.Range("N132").Value = -24783.5

# ---------------------------------------------------------------------------
# WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H126").Value = 3224.9333
$ws.Range("I126").Value = 3141.92
$ws.Range("J126").Value = 3640
$ws.Range("K126").Value = 9425.76
$ws.Range("L126").Value = 10920
$ws.Range("M126").Value = -6955.76
$ws.Range("N126").Value = -15860
